# Update "想去人数" (F column) figures on both the "展览" and "全部类型" sheets
# to match the newly scraped output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 15007
$ws1.Range("F3").Value  = 19034
$ws1.Range("F5").Value  = 137
$ws1.Range("F14").Value = 159
$ws1.Range("F15").Value = 218
$ws1.Range("F17").Value = 1465
$ws1.Range("F22").Value = 7931
$ws1.Range("F23").Value = 989
$ws1.Range("F27").Value = 1246
$ws1.Range("F28").Value = 17
$ws1.Range("F29").Value = 6053
$ws1.Range("F30").Value = 115
$ws1.Range("F34").Value = 283
$ws1.Range("F35").Value = 5426
$ws1.Range("F36").Value = 236
$ws1.Range("F37").Value = 11

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 15007
$ws4.Range("F3").Value  = 19034
$ws4.Range("F5").Value  = 137
$ws4.Range("F14").Value = 159
$ws4.Range("F15").Value = 218
$ws4.Range("F17").Value = 1465
$ws4.Range("F23").Value = 7931
$ws4.Range("F24").Value = 989
$ws4.Range("F28").Value = 1246
$ws4.Range("F29").Value = 17
$ws4.Range("F32").Value = 6053
$ws4.Range("F33").Value = 115
$ws4.Range("F37").Value = 283
$ws4.Range("F38").Value = 5426
$ws4.Range("F39").Value = 237
$ws4.Range("F40").Value = 11
